# Refresh the cryptocurrency market data table (price + 1h volume change).
# Mirrors automated GitHub Actions commit:
#   "Updated cryptos list on Wed Oct  2 20:58:48 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.898.09'
$ws.Range('E2').Value = '  +0.09%  '
# Row 3
$ws.Range('D3').Value = '2.386.47'
$ws.Range('E3').Value = '  -2.89%  '
# Row 4
$ws.Range('E4').Value = '  +0.11%  '
# Row 5
$ws.Range('D5').Value = "'544.47"
$ws.Range('E5').Value = '  -0.07%  '
# Row 6
$ws.Range('D6').Value = "'140.97"
$ws.Range('E6').Value = '  -3.33%  '
# Row 7
$ws.Range('E7').Value = '  +0.01%  '
# Row 8
$ws.Range('D8').Value = "'0.575"
$ws.Range('E8').Value = '  -5.66%  '
# Row 9
$ws.Range('D9').Value = '2.389.65'
$ws.Range('E9').Value = '  -2.71%  '
# Row 10
$ws.Range('E10').Value = '  -1.15%  '
# Row 11
$ws.Range('E11').Value = '  +0.65%  '
# Row 12
$ws.Range('D12').Value = "'5.37"
$ws.Range('E12').Value = '  +0.50%  '
# Row 13
$ws.Range('D13').Value = "'0.344"
$ws.Range('E13').Value = '  -2.86%  '
# Row 14
$ws.Range('D14').Value = "'25.50"
$ws.Range('E14').Value = '  -1.65%  '
# Row 15
$ws.Range('D15').Value = '2.815.15'
$ws.Range('E15').Value = '  -2.60%  '
# Row 16
$ws.Range('E16').Value = '  +0.76%  '
# Row 17
$ws.Range('D17').Value = '60.539.03'
$ws.Range('E17').Value = '  -0.31%  '
# Row 18
$ws.Range('D18').Value = '2.382.76'
$ws.Range('E18').Value = '  -2.81%  '
# Row 19
$ws.Range('D19').Value = "'10.64"
$ws.Range('E19').Value = '  -3.97%  '
# Row 20
$ws.Range('D20').Value = "'4.12"
$ws.Range('E20').Value = '  -1.73%  '
# Row 21
$ws.Range('D21').Value = "'318.16"
$ws.Range('E21').Value = '  -0.08%  '
# Row 22
$ws.Range('D22').Value = "'6.71"
$ws.Range('E22').Value = '  -3.58%  '
# Row 23
$ws.Range('E23').Value = '  -0.09%  '
# Row 24
$ws.Range('E24').Value = '  +3.70%  '
# Row 25
$ws.Range('D25').Value = "'63.27"
$ws.Range('E25').Value = '  +0.29%  '
# Row 26
$ws.Range('D26').Value = "'0.998"
$ws.Range('E26').Value = '  -0.12%  '
# Row 27
$ws.Range('D27').Value = '2.503.73'
$ws.Range('E27').Value = '  -3.05%  '
# Row 28
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0937'
$ws.Range('E28').Value = '  -4.07%  '
# Row 29
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = "'7.84"
$ws.Range('E29').Value = '  +2.30%  '
# Row 30
$ws.Range('D30').Value = "'523.47"
$ws.Range('E30').Value = '  -1.92%  '
# Row 31
$ws.Range('E31').Value = '  -4.13%  '
# Row 32
$ws.Range('D32').Value = "'8.03"
$ws.Range('E32').Value = '  -3.62%  '
# Row 33
$ws.Range('E33').Value = '  -1.51%  '
# Row 34
$ws.Range('E34').Value = '  -2.74%  '
# Row 35
$ws.Range('E35').Value = '  +0.14%  '
# Row 36
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  +0.16%  '
# Row 37
$ws.Range('D37').Value = "'5.48"
$ws.Range('E37').Value = '  -6.11%  '
# Row 38
$ws.Range('D38').Value = "'4.67"
$ws.Range('E38').Value = '  -4.09%  '
# Row 39
$ws.Range('E39').Value = '  +0.52%  '
# Row 40
$ws.Range('D40').Value = "'18.12"
$ws.Range('E40').Value = '  -1.22%  '
# Row 41
$ws.Range('E41').Value = '  +1.90%  '
# Row 42
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = "'138.37"
$ws.Range('E42').Value = '  -3.86%  '
# Row 43
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  +0.19%  '
# Row 44
$ws.Range('D44').Value = "'40.23"
$ws.Range('E44').Value = '  +0.66%  '
# Row 45
$ws.Range('D45').Value = "'2.24"
$ws.Range('E45').Value = '  -3.77%  '
# Row 46
$ws.Range('D46').Value = "'140.25"
$ws.Range('E46').Value = '  -4.33%  '
# Row 47
$ws.Range('D47').Value = "'3.55"
$ws.Range('E47').Value = '  -0.58%  '
# Row 48
$ws.Range('D48').Value = "'20.43"
$ws.Range('E48').Value = '  -2.38%  '
# Row 49
$ws.Range('D49').Value = "'0.0516"
$ws.Range('E49').Value = '  -2.94%  '
# Row 50
$ws.Range('D50').Value = "'0.580"
$ws.Range('E50').Value = '  -0.25%  '
# Row 51
$ws.Range('D51').Value = "'0.0926"
$ws.Range('E51').Value = '  -1.57%  '

